$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.406.01'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.35%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.878.96'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7175'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.23%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '243.78'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.76%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07931'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.55%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.94'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.68%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08135'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.03%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.878.60'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.27%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '95.45'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +4.57%  '

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.07%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.7074'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.41%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.417'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +4.43%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008408'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.16%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '29.406.87'

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '252.52'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.82%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.39'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.41%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.140.13'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.74%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.658'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -1.26%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.12%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1585'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.56%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.070'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.31%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '162.28'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.35%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.90'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.06%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.20%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.416'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.18%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.294'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.04%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.216'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.84%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05320'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.63%  '

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.08%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7577'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.76%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.176'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.700'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.69%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01893'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.94%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.268.89'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.94%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.760'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.95%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.397'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.68%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '111.99'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.65%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9049'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.32%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '74.24'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.72%  '

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.15%  '

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.25%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.035.89'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.52%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.811'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.65%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.5206'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.28%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.513'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.64%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4346'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.17%  '
